# Replaces Avij vars with Avii' vars in month model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the base input parameters (B1:B3) ---
$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 18

# --- Build the new right-hand block (columns G,H,J,K), mirroring A/B/D/E ---

# Row 11 headers/totals (bold, like A11/D11)
$ws.Range("G11").Value = "Rows"
$ws.Range("G11").Font.Bold = $true
$ws.Range("H11").Formula = "=SUM(H12:H17)"

$ws.Range("J11").Value = "NonZero"
$ws.Range("J11").Font.Bold = $true
$ws.Range("K11").Formula = "=SUM(K12:K17)"

# Row 12
$ws.Range("G12").Value = "Limit"
$ws.Range("H12").Formula = "=B2"
$ws.Range("J12").Formula = "=B5"
$ws.Range("K12").Formula = "=H12*J12"

# Row 13
$ws.Range("G13").Value = "Order"
$ws.Range("H13").Formula = "=B2*B5"
$ws.Range("J13").Formula = "=2+2*(B5-1)"
$ws.Range("K13").Formula = "=H13*J13"

# Row 14
$ws.Range("G14").Value = "Resource"
$ws.Range("H14").Formula = "=B1*B3"
$ws.Range("J14").Formula = "=(B2/B1)*B5"
$ws.Range("K14").Formula = "=H14*J14"

# Row 15
$ws.Range("G15").Value = "Duration"
$ws.Range("H15").Formula = "=B2*B5*(B5-1)"
$ws.Range("J15").Value = 3
$ws.Range("K15").Formula = "=H15*J15"

# Row 16
$ws.Range("G16").Value = "Finish"
$ws.Range("H16").Formula = "=B3"
$ws.Range("J16").Value = 1
$ws.Range("K16").Formula = "=H16*J16"

# Row 17 (mirrors A18/"Fixed" row)
$ws.Range("G17").Value = "Fixed"
$ws.Range("H17").Formula = "=B4*B2"
$ws.Range("J17").Formula = "=B5"
$ws.Range("K17").Formula = "=H17*J17"

# Remove the old empty placeholder cell at G18 (no longer present in the
# new layout; the "Fixed" row now lives at G17/H17 above)
$ws.Range("G18").Clear() | Out-Null

# Row 19 (new row, bold like A20/"Columns" header)
$ws.Range("G19").Value = "Columns"
$ws.Range("G19").Font.Bold = $true
$ws.Range("H19").Formula = "=H20+H21+H22+H23"

# Row 20 (replace the old bold empty placeholder with real, unstyled content)
$ws.Range("G20").Clear() | Out-Null
$ws.Range("G20").Value = "s"
$ws.Range("H20").Formula = "=B5"

# Row 21
$ws.Range("G21").Value = "a"
$ws.Range("H21").Formula = "=B2*B5*B5"

# Row 23 (write before row 22 so the new shared-strings are appended in
# the same order as the target workbook: aL=17, aF=18)
$ws.Range("G23").Value = "aL"
$ws.Range("H23").Formula = "=B2*B5"

# Row 22
$ws.Range("G22").Value = "aF"
$ws.Range("H22").Formula = "=B2*B5"

# --- Selection / view state ---
$ws.Range("B5").Select() | Out-Null

$wb.Application.CalculateFull()
